# Apply edits to the "06组项目计划表" workbook (Sheet1 is the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 43 (邱培松): shorten task description.
$ws.Range("B43").Value = "内容:完成数据库物理模型图"

# Row 33 (邱培松): shorten task description.
$ws.Range("B33").Value = "内容:设计数据库逻辑模型"

# Row 44 (黄龙强): update completion percentage.
$ws.Range("C44").Value = 0.5

# Leave the merged range A27:D28 selected, matching the saved view state.
$ws.Activate()
$ws.Range("A27:D28").Select()
